# Update crypto price/volume figures per the Wed Oct 4 21:41:28 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force a literal text value even when it looks numeric (e.g. "213.42"),
    # matching the inline/shared string cells already in the sheet, then drop
    # the temporary text number-format so no stray style sticks to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '27.756.15'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '1.646.08'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue $ws.Range("D5") '213.42'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").Value = '  +3.67%  '
Set-TextValue $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  -0.07%  '
Set-TextValue $ws.Range("D8") '23.14'
$ws.Range("E8").Value = '  -2.10%  '
$ws.Range("E9").Value = '  +0.04%  '
Set-TextValue $ws.Range("D10") '0.0614'
$ws.Range("E10").Value = '  +0.07%  '
Set-TextValue $ws.Range("D11") '0.0891'
$ws.Range("E11").Value = '  +1.83%  '
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("D13").Value = '1.657.94'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("E15").Value = '  -1.03%  '
Set-TextValue $ws.Range("D16") '64.38'
$ws.Range("E16").Value = '  -1.76%  '
$ws.Range("D17").Value = '27.715.43'
$ws.Range("E17").Value = '  +1.22%  '
Set-TextValue $ws.Range("D18") '231.92'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = '0.0₃0727'
$ws.Range("E19").Value = '  +0.13%  '
Set-TextValue $ws.Range("D20") '7.68'
$ws.Range("E20").Value = '  +3.03%  '
$ws.Range("E21").Value = '  +0.05%  '
Set-TextValue $ws.Range("D22") '4.33'
$ws.Range("E22").Value = '  -0.72%  '
Set-TextValue $ws.Range("D23") '10.11'
$ws.Range("E23").Value = '  +7.98%  '
Set-TextValue $ws.Range("D24") '1.97'
$ws.Range("E24").Value = '  -2.97%  '
Set-TextValue $ws.Range("D25") '149.61'
$ws.Range("E25").Value = '  +1.25%  '
Set-TextValue $ws.Range("D26") '6.99'
$ws.Range("E26").Value = '  -1.57%  '
$ws.Range("E27").Value = '  +0.94%  '
Set-TextValue $ws.Range("D28") '15.70'
$ws.Range("E28").Value = '  -1.03%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  -0.17%  '
Set-TextValue $ws.Range("D31") '0.0487'
$ws.Range("E31").Value = '  -2.12%  '
Set-TextValue $ws.Range("D32") '3.31'
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("E33").Value = '  +1.59%  '
$ws.Range("D34").Value = '1.445.72'
$ws.Range("E34").Value = '  +1.47%  '
Set-TextValue $ws.Range("D35") '1.60'
$ws.Range("E35").Value = '  +2.19%  '
$ws.Range("E36").Value = '  -1.17%  '
Set-TextValue $ws.Range("D37") '0.573'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  -2.49%  '
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("E40").Value = '  +12.34%  '
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("E42").Value = '  +0.06%  '
Set-TextValue $ws.Range("D43") '5.66'
$ws.Range("E43").Value = '  +2.79%  '
$ws.Range("E44").Value = '  -0.43%  '
$ws.Range("E45").Value = '  +2.03%  '
Set-TextValue $ws.Range("D46") '65.86'
$ws.Range("E46").Value = '  +1.45%  '
$ws.Range("D47").Value = '1.787.88'
$ws.Range("E47").Value = '  -0.58%  '
Set-TextValue $ws.Range("D48") '1.71'
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("E49").Value = '  +1.59%  '
Set-TextValue $ws.Range("D50") '85.96'
$ws.Range("E50").Value = '  -2.31%  '
Set-TextValue $ws.Range("D51") '0.0993'
$ws.Range("E51").Value = '  -1.67%  '
